# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (it already has the right column layout
#    and cell styles) to create the new "2022-Q1" sheet, positioned right
#    before the "总计" sheet.
# 2. Trim / overwrite its rows with the 2022-Q1 holdings data.
# 3. Insert a new top data row into the "总计" sheet summarizing the
#    2022-Q1 quarter, shifting the existing rows down.
# 4. Restore the originally-selected tab.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$totalBeforeInsert = $wb.Worksheets.Item("总计")

# --- Step 1: create the new sheet by copying "2021-Q4" (keeps formatting) ---
# NOTE: worksheet references here behave like live index pointers, so as
# soon as a sheet is inserted/removed, any previously-captured reference
# whose position shifted (like $totalBeforeInsert) can silently start
# pointing at a different sheet. Always re-fetch sheets *by name* right
# after a sheet-count-changing operation instead of reusing old handles.
$q4.Copy($totalBeforeInsert)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# the copied sheet has 8 rows (header + 7 data rows); the new sheet only
# needs 7 rows (header + 6 data rows), so drop the extra trailing row
$q1.Rows.Item(8).Delete()

# --- Step 2: fill in the 2022-Q1 holdings data ---
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

$data = @(
    @("516970", "广发中证基建工程交易型开放式指数证券投资基金", "59.01", "99.38", "5.45", "3.2160", 6),
    @("165525", "信诚中证基建工程指数（LOF）",                   "17.06", "94.00", "5.16", "0.8803", 6),
    @("516950", "银华中证基建交易型开放式指数证券投资基金",       "10.41", "97.55", "4.05", "0.4216", 8),
    @("510081", "长盛动态精选混合",                             "3.15",  "60.76", "4.27", "0.1345", 3),
    @("006478", "长盛多因子策略优选股票",                       "0.51",  "84.41", "4.67", "0.0238", 5),
    @("003238", "新华外延增长主题灵活配置混合",                 "0.42",  "87.12", "2.47", "0.0104", 7)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $item = $data[$i]

    $q1.Cells.Item($row,1).Value = $i

    # H (仓位排名) is a genuine number; set it first since, being left at
    # the sheet's default style, it also doubles as a "clean" formatting
    # source for the text cells below
    $q1.Cells.Item($row,8).Value = $item[6]

    # columns B, D, E, F, G hold text values (fund codes with leading
    # zeros / numeric-looking figures such as "94.00" or "3.2160" whose
    # formatting must be preserved exactly). Temporarily force the cell to
    # Text format so the literal string is kept instead of being parsed
    # into a number, then restore the cell's original (default) format by
    # pasting the formatting from the untouched H cell in the same row so
    # no stray style is left behind.
    $q1.Cells.Item($row,2).NumberFormat = "@"
    $q1.Cells.Item($row,2).Value = $item[0]
    $q1.Cells.Item($row,8).Copy()
    $q1.Cells.Item($row,2).PasteSpecial(-4122)

    $q1.Cells.Item($row,3).Value = $item[1]

    $q1.Cells.Item($row,4).NumberFormat = "@"
    $q1.Cells.Item($row,4).Value = $item[2]
    $q1.Cells.Item($row,8).Copy()
    $q1.Cells.Item($row,4).PasteSpecial(-4122)

    $q1.Cells.Item($row,5).NumberFormat = "@"
    $q1.Cells.Item($row,5).Value = $item[3]
    $q1.Cells.Item($row,8).Copy()
    $q1.Cells.Item($row,5).PasteSpecial(-4122)

    $q1.Cells.Item($row,6).NumberFormat = "@"
    $q1.Cells.Item($row,6).Value = $item[4]
    $q1.Cells.Item($row,8).Copy()
    $q1.Cells.Item($row,6).PasteSpecial(-4122)

    $q1.Cells.Item($row,7).NumberFormat = "@"
    $q1.Cells.Item($row,7).Value = $item[5]
    $q1.Cells.Item($row,8).Copy()
    $q1.Cells.Item($row,7).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# --- Step 3: update the "总计" sheet with the new 2022-Q1 summary row ---
# re-fetch "总计" by name now that the sheet collection has changed
$total = $wb.Worksheets.Item("总计")

# copy the number-row style (bold/border, column A) down onto the new
# last row before the values get shifted
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(4,1).PasteSpecial(-4122)

# read the existing rows 2 and 3 before overwriting them
$oldB2 = $total.Cells.Item(2,2).Value()
$oldC2 = $total.Cells.Item(2,3).Value()
$oldD2 = $total.Cells.Item(2,4).Value()
$oldB3 = $total.Cells.Item(3,2).Value()
$oldC3 = $total.Cells.Item(3,3).Value()
$oldD3 = $total.Cells.Item(3,4).Value()

# row 4 = old row 3 ("2021-Q3")
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = $oldB3
$total.Cells.Item(4,3).Value = $oldC3
$total.Cells.Item(4,4).Value = $oldD3

# row 3 = old row 2 ("2021-Q4")
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = $oldB2
$total.Cells.Item(3,3).Value = $oldC2
$total.Cells.Item(3,4).Value = $oldD2

# row 2 = new "2022-Q1" summary
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,4).Value = 4.69

# --- Step 4: restore the originally active tab ---
$wb.Worksheets.Item("2021-Q3").Activate()
